# Regen sval data to filter save games
# Updates numeric columns B:G for rows 2-8 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    3 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 0, 18.67282528286833)
    4 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 3.754798637575387)
    5 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 0, 14.05633640148523)
    6 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    7 = @(0.1190320826869504, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 0, 4.457851494814137)
    8 = @(0.2917716402565462, 0.00006240767534437808, 0.1494219747398047, 0.4942365360607697, 0, 0.935492558732465)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
